# RPA datasets push 2024-05-25
# Insert two new IPO book-building rows for "에이치브이엠(구.한국진공야금)" and
# "이노스페이스" right after row 3, then remove the two stale rows at the
# bottom of the table that used to hold those same companies (their data has
# been refreshed / moved to the new rows near the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows (old row4 -> row6, ... old row21 -> row23).
$ws.Rows("4:5").Insert()

$ws.Cells.Item(4, 1).Value = "에이치브이엠(구.한국진공야금)"
$ws.Cells.Item(4, 2).Value = "2024.06.11~06.17"
$ws.Cells.Item(4, 3).Value = "11,000~14,200"
$ws.Cells.Item(4, 4).Value = "-"
$e4 = $ws.Cells.Item(4, 5)
$e4.NumberFormat = "@"
$e4.Value = "26400"
$e4.Style = "Normal"
$ws.Cells.Item(4, 6).Value = "NH투자증권"

$ws.Cells.Item(5, 1).Value = "이노스페이스"
$ws.Cells.Item(5, 2).Value = "2024.06.11~06.17"
$ws.Cells.Item(5, 3).Value = "36,400~43,300"
$ws.Cells.Item(5, 4).Value = "-"
$e5 = $ws.Cells.Item(5, 5)
$e5.NumberFormat = "@"
$e5.Value = "48412"
$e5.Style = "Normal"
$ws.Cells.Item(5, 6).Value = "미래에셋증권,신한투자증권"

# The same two companies' old rows (now at 21 and 22 after the insert above)
# are now redundant -- delete them so the table stays at 20 data rows.
$ws.Rows("21:22").Delete()
